$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force column A to Text format for the new rows so date-like strings
# ("06-08-2021", etc.) are stored as literal text, matching the source data,
# instead of being auto-converted to date serial numbers.
$ws.Range("A141:A162").NumberFormat = "@"

$ws.Range("A141").Value = "06-08-2021"
$ws.Range("B141").Value = 40
$ws.Range("C141").Value = 95
$ws.Range("D141").Value = 40
$ws.Range("E141").Value = 787

$ws.Range("A142").Value = "09-08-2021"
$ws.Range("B142").Value = 40
$ws.Range("C142").Value = 125
$ws.Range("D142").Value = 30
$ws.Range("E142").Value = 788

$ws.Range("A143").Value = "10-08-2021"
$ws.Range("B143").Value = 50
$ws.Range("C143").Value = 205
$ws.Range("D143").Value = 50
$ws.Range("E143").Value = 780

$ws.Range("A144").Value = "11-08-2021"
$ws.Range("B144").Value = 40
$ws.Range("C144").Value = 156
$ws.Range("D144").Value = 40
$ws.Range("E144").Value = 774

$ws.Range("A145").Value = "12-08-2021"
$ws.Range("B145").Value = 40
$ws.Range("C145").Value = 141
$ws.Range("D145").Value = 40
$ws.Range("E145").Value = 775

$ws.Range("A146").Value = "13-08-2021"
$ws.Range("B146").Value = 40
$ws.Range("C146").Value = 141
$ws.Range("D146").Value = 40
$ws.Range("E146").Value = 775

$ws.Range("A147").Value = "16-08-2021"
$ws.Range("B147").Value = 40
$ws.Range("C147").Value = 137
$ws.Range("D147").Value = 40
$ws.Range("E147").Value = 789

$ws.Range("A148").Value = "17-08-2021"
$ws.Range("B148").Value = 40
$ws.Range("C148").Value = 160
$ws.Range("D148").Value = 40
$ws.Range("E148").Value = 791

$ws.Range("A149").Value = "18-08-2021"
$ws.Range("B149").Value = 40
$ws.Range("C149").Value = 160
$ws.Range("D149").Value = 40
$ws.Range("E149").Value = 787

$ws.Range("A150").Value = "19-08-2021"
$ws.Range("B150").Value = 40
$ws.Range("C150").Value = 170
$ws.Range("D150").Value = 40
$ws.Range("E150").Value = 791

$ws.Range("A151").Value = "20-08-2021"
$ws.Range("B151").Value = 40
$ws.Range("C151").Value = 121
$ws.Range("D151").Value = 40
$ws.Range("E151").Value = 787

$ws.Range("A152").Value = "23-08-2021"
$ws.Range("B152").Value = 40
$ws.Range("C152").Value = 95
$ws.Range("D152").Value = 35
$ws.Range("E152").Value = 784

$ws.Range("A153").Value = "24-08-2021"
$ws.Range("B153").Value = 45
$ws.Range("C153").Value = 147
$ws.Range("D153").Value = 45
$ws.Range("E153").Value = 782

$ws.Range("A154").Value = "25-08-2021"
$ws.Range("B154").Value = 40
$ws.Range("C154").Value = 179
$ws.Range("D154").Value = 40
$ws.Range("E154").Value = 784

$ws.Range("A155").Value = "26-08-2021"
$ws.Range("B155").Value = 40
$ws.Range("C155").Value = 144
$ws.Range("D155").Value = 40
$ws.Range("E155").Value = 786

$ws.Range("A156").Value = "27-08-2021"
$ws.Range("B156").Value = 40
$ws.Range("C156").Value = 151
$ws.Range("D156").Value = 40
$ws.Range("E156").Value = 784

$ws.Range("A157").Value = "30-08-2021"
$ws.Range("B157").Value = 40
$ws.Range("C157").Value = 127
$ws.Range("D157").Value = 40
$ws.Range("E157").Value = 780

$ws.Range("A158").Value = "31-08-2021"
$ws.Range("B158").Value = 40
$ws.Range("C158").Value = 140
$ws.Range("D158").Value = 40
$ws.Range("E158").Value = 776

$ws.Range("A159").Value = "01-09-2021"
$ws.Range("B159").Value = 40
$ws.Range("C159").Value = 162
$ws.Range("D159").Value = 40
$ws.Range("E159").Value = 769

$ws.Range("A160").Value = "02-09-2021"
$ws.Range("B160").Value = 40
$ws.Range("C160").Value = 153
$ws.Range("D160").Value = 40
$ws.Range("E160").Value = 767

$ws.Range("A161").Value = "03-09-2021"
$ws.Range("B161").Value = 40
$ws.Range("C161").Value = 148
$ws.Range("D161").Value = 40
$ws.Range("E161").Value = 767

$ws.Range("A162").Value = "06-09-2021"
$ws.Range("B162").Value = 40
$ws.Range("C162").Value = 169
$ws.Range("D162").Value = 40
$ws.Range("E162").Value = 772

# Restore the default ("Normal") cell style on column A for these rows so the
# cells do not carry an explicit style index (matches the rest of the sheet,
# where data cells have no "s" attribute).
$ws.Range("A141:A162").Style = "Normal"
